# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 2311.3  # H38: 3391.8572 -> 2311.3
$ws.Cells.Item(38, 9).Value = 1790.8889  # I38: 2549.8 -> 1790.8889
$ws.Cells.Item(38, 10).Value = 6995  # J38: 5497 -> 6995
$ws.Cells.Item(38, 11).Value = 5372.6667  # K38: 7649.400000000001 -> 5372.6667
$ws.Cells.Item(38, 12).Value = 20985  # L38: 16491 -> 20985
$ws.Cells.Item(38, 13).Value = -5000.6667  # M38: -7277.400000000001 -> -5000.6667
$ws.Cells.Item(38, 14).Value = -21729  # N38: -17235 -> -21729
$ws.Cells.Item(40, 8).Value = 10018.154  # H40: 9695.799999999999 -> 10018.154
$ws.Cells.Item(40, 9).Value = 8250.5  # I40: 10001 -> 8250.5
$ws.Cells.Item(40, 10).Value = 10803.777  # J40: 9619.5 -> 10803.777
$ws.Cells.Item(40, 11).Value = 8250.5  # K40: 10001 -> 8250.5
$ws.Cells.Item(40, 12).Value = 10803.777  # L40: 9619.5 -> 10803.777
$ws.Cells.Item(40, 13).Value = -8075.5  # M40: -9826 -> -8075.5
$ws.Cells.Item(40, 14).Value = -11153.777  # N40: -9969.5 -> -11153.777
$ws.Cells.Item(41, 9).Value = 619.8570999999999  # I41: 620.8570999999999 -> 619.8570999999999
$ws.Cells.Item(41, 11).Value = 619.8570999999999  # K41: 620.8570999999999 -> 619.8570999999999
$ws.Cells.Item(41, 13).Value = -179.8570999999999  # M41: -180.8570999999999 -> -179.8570999999999
$ws.Cells.Item(58, 8).Value = 4303.1333  # H58: 4046.5386 -> 4303.1333
$ws.Cells.Item(58, 9).Value = 519.1111  # I58: 337.25 -> 519.1111
$ws.Cells.Item(58, 10).Value = 9979.166999999999  # J58: 9981.4 -> 9979.166999999999
$ws.Cells.Item(58, 11).Value = 1557.3333  # K58: 1011.75 -> 1557.3333
$ws.Cells.Item(58, 12).Value = 29937.501  # L58: 29944.2 -> 29937.501
$ws.Cells.Item(58, 13).Value = -1407.3333  # M58: -861.75 -> -1407.3333
$ws.Cells.Item(58, 14).Value = -30237.501  # N58: -30244.2 -> -30237.501
$ws.Cells.Item(64, 8).Value = 10224.5  # H64: 6197.909 -> 10224.5
$ws.Cells.Item(64, 9).Value = 3898  # I64: 3897.125 -> 3898
$ws.Cells.Item(64, 11).Value = 3898  # K64: 3897.125 -> 3898
$ws.Cells.Item(64, 13).Value = -3650  # M64: -3649.125 -> -3650
$ws.Cells.Item(67, 8).Value = 10224.5  # H67: 6197.909 -> 10224.5
$ws.Cells.Item(67, 9).Value = 3898  # I67: 3897.125 -> 3898
$ws.Cells.Item(67, 11).Value = 3898  # K67: 3897.125 -> 3898
$ws.Cells.Item(67, 13).Value = -3040  # M67: -3039.125 -> -3040
$ws.Cells.Item(80, 8).Value = 4396.161  # H80: 4656.2856 -> 4396.161
$ws.Cells.Item(80, 10).Value = 1560.1333  # J80: 1458.0834 -> 1560.1333
$ws.Cells.Item(80, 12).Value = 4680.3999  # L80: 4374.2502 -> 4680.3999
$ws.Cells.Item(80, 14).Value = -6676.3999  # N80: -6370.2502 -> -6676.3999
$ws.Cells.Item(82, 8).Value = 900.1818  # H82: 1046.5 -> 900.1818
$ws.Cells.Item(82, 9).Value = 900.1818  # I82: 1046.5 -> 900.1818
$ws.Cells.Item(82, 11).Value = 2700.5454  # K82: 3139.5 -> 2700.5454
$ws.Cells.Item(82, 13).Value = -2294.5454  # M82: -2733.5 -> -2294.5454
$ws.Cells.Item(83, 8).Value = 4396.161  # H83: 4656.2856 -> 4396.161
$ws.Cells.Item(83, 10).Value = 1560.1333  # J83: 1458.0834 -> 1560.1333
$ws.Cells.Item(83, 12).Value = 14041.1997  # L83: 13122.7506 -> 14041.1997
$ws.Cells.Item(83, 14).Value = -24025.1997  # N83: -23106.7506 -> -24025.1997
$ws.Cells.Item(85, 8).Value = 900.1818  # H85: 1046.5 -> 900.1818
$ws.Cells.Item(85, 9).Value = 900.1818  # I85: 1046.5 -> 900.1818
$ws.Cells.Item(85, 11).Value = 2700.5454  # K85: 3139.5 -> 2700.5454
$ws.Cells.Item(85, 13).Value = -1296.5454  # M85: -1735.5 -> -1296.5454
$ws.Cells.Item(87, 8).Value = 71953.84  # H87: 67943.75 -> 71953.84
$ws.Cells.Item(87, 9).Value = 42000  # I87: 40000 -> 42000
$ws.Cells.Item(87, 10).Value = 74450  # J87: 74392.30499999999 -> 74450
$ws.Cells.Item(87, 11).Value = 42000  # K87: 40000 -> 42000
$ws.Cells.Item(87, 12).Value = 74450  # L87: 74392.30499999999 -> 74450
$ws.Cells.Item(87, 13).Value = -40752  # M87: -38752 -> -40752
$ws.Cells.Item(87, 14).Value = -76946  # N87: -76888.30499999999 -> -76946
$ws.Cells.Item(88, 8).Value = 4698.7896  # H88: 6060.769 -> 4698.7896
$ws.Cells.Item(88, 9).Value = 5759.2  # I88: 8333.333000000001 -> 5759.2
$ws.Cells.Item(88, 10).Value = 4320.0713  # J88: 5379 -> 4320.0713
$ws.Cells.Item(88, 11).Value = 5759.2  # K88: 8333.333000000001 -> 5759.2
$ws.Cells.Item(88, 12).Value = 4320.0713  # L88: 5379 -> 4320.0713
$ws.Cells.Item(88, 13).Value = -5353.2  # M88: -7927.333000000001 -> -5353.2
$ws.Cells.Item(88, 14).Value = -5132.0713  # N88: -6191 -> -5132.0713
$ws.Cells.Item(90, 8).Value = 71953.84  # H90: 67943.75 -> 71953.84
$ws.Cells.Item(90, 9).Value = 42000  # I90: 40000 -> 42000
$ws.Cells.Item(90, 10).Value = 74450  # J90: 74392.30499999999 -> 74450
$ws.Cells.Item(90, 11).Value = 126000  # K90: 120000 -> 126000
$ws.Cells.Item(90, 12).Value = 223350  # L90: 223176.915 -> 223350
$ws.Cells.Item(90, 13).Value = -119760  # M90: -113760 -> -119760
$ws.Cells.Item(90, 14).Value = -235830  # N90: -235656.915 -> -235830
$ws.Cells.Item(91, 8).Value = 4698.7896  # H91: 6060.769 -> 4698.7896
$ws.Cells.Item(91, 9).Value = 5759.2  # I91: 8333.333000000001 -> 5759.2
$ws.Cells.Item(91, 10).Value = 4320.0713  # J91: 5379 -> 4320.0713
$ws.Cells.Item(91, 11).Value = 5759.2  # K91: 8333.333000000001 -> 5759.2
$ws.Cells.Item(91, 12).Value = 4320.0713  # L91: 5379 -> 4320.0713
$ws.Cells.Item(91, 13).Value = -4355.2  # M91: -6929.333000000001 -> -4355.2
$ws.Cells.Item(91, 14).Value = -7128.0713  # N91: -8187 -> -7128.0713
$ws.Cells.Item(99, 8).Value = 416.6  # H99: 354.75 -> 416.6
$ws.Cells.Item(99, 9).Value = 416.6  # I99: 354.75 -> 416.6
$ws.Cells.Item(99, 11).Value = 1249.8  # K99: 1064.25 -> 1249.8
$ws.Cells.Item(99, 13).Value = 248.1999999999998  # M99: 433.75 -> 248.1999999999998
$ws.Cells.Item(112, 8).Value = 3883.647  # H112: 3866.3 -> 3883.647
$ws.Cells.Item(112, 10).Value = 4082.625  # J112: 4032.9473 -> 4082.625
$ws.Cells.Item(112, 12).Value = 12247.875  # L112: 12098.8419 -> 12247.875
$ws.Cells.Item(112, 14).Value = -14463.875  # N112: -14314.8419 -> -14463.875
$ws.Cells.Item(115, 8).Value = 421.25  # H115: 466.66666 -> 421.25
$ws.Cells.Item(115, 9).Value = 228.33333  # I115: 200 -> 228.33333
$ws.Cells.Item(115, 11).Value = 684.99999  # K115: 600 -> 684.99999
$ws.Cells.Item(115, 13).Value = 882.00001  # M115: 967 -> 882.00001
$ws.Cells.Item(118, 8).Value = 27484  # H118: 32166.666 -> 27484
$ws.Cells.Item(118, 9).Value = 38340  # I118: 47875 -> 38340
$ws.Cells.Item(118, 10).Value = 344  # J118: 750 -> 344
$ws.Cells.Item(118, 11).Value = 115020  # K118: 143625 -> 115020
$ws.Cells.Item(118, 12).Value = 1032  # L118: 2250 -> 1032
$ws.Cells.Item(118, 13).Value = -113363  # M118: -141968 -> -113363
$ws.Cells.Item(118, 14).Value = -4346  # N118: -5564 -> -4346
$ws.Cells.Item(127, 8).Value = 2067.6667  # H127: 1887.8 -> 2067.6667
$ws.Cells.Item(127, 9).Value = 1851.8334  # I127: 1625.7142 -> 1851.8334
$ws.Cells.Item(127, 11).Value = 5555.5002  # K127: 4877.142599999999 -> 5555.5002
$ws.Cells.Item(127, 13).Value = -595.5002000000004  # M127: 82.85740000000078 -> -595.5002000000004
$ws.Cells.Item(129, 8).Value = 2348.2856  # H129: 709.875 -> 2348.2856
$ws.Cells.Item(129, 9).Value = 1146  # I129: 709.875 -> 1146
$ws.Cells.Item(129, 10).Value = 3250  # J129: 0 -> 3250
$ws.Cells.Item(129, 11).Value = 3438  # K129: 2129.625 -> 3438
$ws.Cells.Item(129, 12).Value = 9750  # L129: 0 -> 9750
$ws.Cells.Item(129, 13).Value = 1562  # M129: 2870.375 -> 1562
$ws.Cells.Item(129, 14).Value = -19750  # N129: None -> -19750
$ws.Cells.Item(138, 8).Value = 3882.2273  # H138: 3982.8572 -> 3882.2273
$ws.Cells.Item(138, 9).Value = 820.0833  # I138: 840.56525 -> 820.0833
$ws.Cells.Item(138, 10).Value = 5632.024  # J138: 5789.675 -> 5632.024
$ws.Cells.Item(138, 11).Value = 2460.2499  # K138: 2521.69575 -> 2460.2499
$ws.Cells.Item(138, 12).Value = 16896.072  # L138: 17369.025 -> 16896.072
$ws.Cells.Item(138, 13).Value = 2679.7501  # M138: 2618.30425 -> 2679.7501
$ws.Cells.Item(138, 14).Value = -27176.072  # N138: -27649.025 -> -27176.072

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 894.8333  # H5: 898.1667 -> 894.8333
$ws.Cells.Item(5, 9).Value = 94  # I5: 98 -> 94
$ws.Cells.Item(5, 11).Value = 94  # K5: 98 -> 94
$ws.Cells.Item(5, 13).Value = 18  # M5: 14 -> 18
$ws.Cells.Item(39, 8).Value = 26812.5  # H39: 33758.668 -> 26812.5
$ws.Cells.Item(39, 9).Value = 9899  # I39: 29999.75 -> 9899
$ws.Cells.Item(39, 10).Value = 77553  # J39: 41276.5 -> 77553
$ws.Cells.Item(39, 11).Value = 9899  # K39: 29999.75 -> 9899
$ws.Cells.Item(39, 12).Value = 77553  # L39: 41276.5 -> 77553
$ws.Cells.Item(39, 13).Value = -9379  # M39: -29479.75 -> -9379
$ws.Cells.Item(39, 14).Value = -78593  # N39: -42316.5 -> -78593

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 894.8333  # H4: 898.1667 -> 894.8333
$ws.Cells.Item(4, 9).Value = 94  # I4: 98 -> 94
$ws.Cells.Item(4, 11).Value = 94  # K4: 98 -> 94
$ws.Cells.Item(4, 13).Value = 21  # M4: 17 -> 21
$ws.Cells.Item(14, 8).Value = 0  # H14: 9 -> 0
$ws.Cells.Item(14, 10).Value = 0  # J14: 9 -> 0
$ws.Cells.Item(14, 12).ClearContents()  # L14: 9 -> (removed)
$ws.Cells.Item(14, 14).Value = 0  # N14: -353 -> 0
$ws.Cells.Item(134, 8).Value = 3730.4348  # H134: 4007.8572 -> 3730.4348
$ws.Cells.Item(134, 9).Value = 2942.5  # I134: 3068.5293 -> 2942.5
$ws.Cells.Item(134, 10).Value = 6567  # J134: 8000 -> 6567
$ws.Cells.Item(134, 11).Value = 8827.5  # K134: 9205.5879 -> 8827.5
$ws.Cells.Item(134, 12).Value = 19701  # L134: 24000 -> 19701
$ws.Cells.Item(134, 13).Value = -6292.5  # M134: -6670.5879 -> -6292.5
$ws.Cells.Item(134, 14).Value = -24771  # N134: -29070 -> -24771

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3247.1  # H62: 3341.2222 -> 3247.1
$ws.Cells.Item(62, 10).Value = 3331.7144  # J62: 3487 -> 3331.7144
$ws.Cells.Item(62, 12).Value = 3331.7144  # L62: 3487 -> 3331.7144
$ws.Cells.Item(62, 14).Value = -4579.7144  # N62: -4735 -> -4579.7144
$ws.Cells.Item(65, 8).Value = 3247.1  # H65: 3341.2222 -> 3247.1
$ws.Cells.Item(65, 10).Value = 3331.7144  # J65: 3487 -> 3331.7144
$ws.Cells.Item(65, 12).Value = 16658.572  # L65: 17435 -> 16658.572
$ws.Cells.Item(65, 14).Value = -22898.572  # N65: -23675 -> -22898.572
$ws.Cells.Item(68, 8).Value = 94800  # H68: 113133.336 -> 94800
$ws.Cells.Item(68, 10).Value = 94800  # J68: 113133.336 -> 94800
$ws.Cells.Item(68, 12).Value = 94800  # L68: 113133.336 -> 94800
$ws.Cells.Item(68, 14).Value = -96298  # N68: -114631.336 -> -96298
$ws.Cells.Item(71, 8).Value = 94800  # H71: 113133.336 -> 94800
$ws.Cells.Item(71, 10).Value = 94800  # J71: 113133.336 -> 94800
$ws.Cells.Item(71, 12).Value = 284400  # L71: 339400.008 -> 284400
$ws.Cells.Item(71, 14).Value = -291888  # N71: -346888.008 -> -291888

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 336.33334  # H17: 693.3333 -> 336.33334
$ws.Cells.Item(17, 9).Value = 4.5  # I17: 0 -> 4.5
$ws.Cells.Item(17, 10).Value = 1000  # J17: 693.3333 -> 1000
$ws.Cells.Item(17, 11).Value = 13.5  # K17: 0 -> 13.5
$ws.Cells.Item(17, 12).Value = 3000  # L17: 2079.9999 -> 3000
$ws.Cells.Item(17, 13).Value = 155.5  # M17: None -> 155.5
$ws.Cells.Item(17, 14).Value = -3338  # N17: -2417.9999 -> -3338
$ws.Cells.Item(34, 8).Value = 65043.61  # H34: 64834.61 -> 65043.61
$ws.Cells.Item(34, 10).Value = 90036.53999999999  # J34: 89747.16 -> 90036.53999999999
$ws.Cells.Item(34, 12).Value = 270109.62  # L34: 269241.48 -> 270109.62
$ws.Cells.Item(34, 14).Value = -270277.62  # N34: -269409.48 -> -270277.62
$ws.Cells.Item(39, 8).Value = 7400  # H39: 7438 -> 7400
$ws.Cells.Item(39, 10).Value = 20000  # J39: 11663.333 -> 20000
$ws.Cells.Item(39, 12).Value = 60000  # L39: 34989.999 -> 60000
$ws.Cells.Item(39, 14).Value = -60588  # N39: -35577.999 -> -60588
$ws.Cells.Item(55, 8).Value = 18000  # H55: 13723.75 -> 18000
$ws.Cells.Item(55, 10).Value = 50000  # J55: 25447.5 -> 50000
$ws.Cells.Item(55, 12).Value = 150000  # L55: 76342.5 -> 150000
$ws.Cells.Item(55, 14).Value = -150354  # N55: -76696.5 -> -150354
$ws.Cells.Item(132, 8).Value = 3191.1724  # H132: 2948.1936 -> 3191.1724
$ws.Cells.Item(132, 9).Value = 971.3570999999999  # I132: 955.5 -> 971.3570999999999
$ws.Cells.Item(132, 10).Value = 5263  # J132: 5707.3076 -> 5263
$ws.Cells.Item(132, 11).Value = 8742.213899999999  # K132: 8599.5 -> 8742.213899999999
$ws.Cells.Item(132, 12).Value = 47367  # L132: 51365.7684 -> 47367
$ws.Cells.Item(132, 13).Value = -6212.213899999999  # M132: -6069.5 -> -6212.213899999999
$ws.Cells.Item(132, 14).Value = -52427  # N132: -56425.7684 -> -52427

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 16400  # H57: 23666.666 -> 16400
$ws.Cells.Item(57, 10).Value = 0  # J57: 60000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 60000 -> 0
$ws.Cells.Item(57, 14).ClearContents()  # N57: -61640 -> (removed)
$ws.Cells.Item(80, 8).Value = 10000003  # H80: 6001802 -> 10000003
$ws.Cells.Item(80, 9).Value = 10000002  # I80: 5002251 -> 10000002
$ws.Cells.Item(80, 11).Value = 10000002  # K80: 5002251 -> 10000002
$ws.Cells.Item(80, 13).Value = -9999004  # M80: -5001253 -> -9999004
$ws.Cells.Item(83, 8).Value = 10000003  # H83: 6001802 -> 10000003
$ws.Cells.Item(83, 9).Value = 10000002  # I83: 5002251 -> 10000002
$ws.Cells.Item(83, 11).Value = 50000010  # K83: 25011255 -> 50000010
$ws.Cells.Item(83, 13).Value = -49995018  # M83: -25006263 -> -49995018
$ws.Cells.Item(132, 8).Value = 438667.66  # H132: 403696.25 -> 438667.66
$ws.Cells.Item(132, 9).Value = 558742  # I132: 503020.3 -> 558742
$ws.Cells.Item(132, 11).Value = 1676226  # K132: 1509060.9 -> 1676226
$ws.Cells.Item(132, 13).Value = -1673696  # M132: -1506530.9 -> -1673696

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1774.5294  # H22: 1717.6111 -> 1774.5294
$ws.Cells.Item(22, 9).Value = 948.1111  # I22: 928.3 -> 948.1111
$ws.Cells.Item(22, 11).Value = 948.1111  # K22: 928.3 -> 948.1111
$ws.Cells.Item(22, 13).Value = -653.1111  # M22: -633.3 -> -653.1111
$ws.Cells.Item(27, 8).Value = 1774.5294  # H27: 1717.6111 -> 1774.5294
$ws.Cells.Item(27, 9).Value = 948.1111  # I27: 928.3 -> 948.1111
$ws.Cells.Item(27, 11).Value = 948.1111  # K27: 928.3 -> 948.1111
$ws.Cells.Item(27, 13).Value = -841.1111  # M27: -821.3 -> -841.1111
$ws.Cells.Item(93, 8).Value = 1285.5  # H93: 1242.7778 -> 1285.5
$ws.Cells.Item(93, 9).Value = 1285.5  # I93: 1242.7778 -> 1285.5
$ws.Cells.Item(93, 11).Value = 1285.5  # K93: 1242.7778 -> 1285.5
$ws.Cells.Item(93, 13).Value = -37.5  # M93: 5.22219999999993 -> -37.5
$ws.Cells.Item(100, 8).Value = 1692.3572  # H100: 1739.7333 -> 1692.3572
$ws.Cells.Item(100, 9).Value = 1668.6923  # I100: 1699.6923 -> 1668.6923
$ws.Cells.Item(100, 11).Value = 1668.6923  # K100: 1699.6923 -> 1668.6923
$ws.Cells.Item(100, 13).Value = -1127.6923  # M100: -1158.6923 -> -1127.6923
$ws.Cells.Item(132, 8).Value = 3928.4138  # H132: 4036.5715 -> 3928.4138
$ws.Cells.Item(132, 9).Value = 1807.4375  # I132: 1867.9333 -> 1807.4375
$ws.Cells.Item(132, 11).Value = 5422.3125  # K132: 5603.7999 -> 5422.3125
$ws.Cells.Item(132, 13).Value = -2892.3125  # M132: -3073.7999 -> -2892.3125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 69374.414  # H62: 56696.617 -> 69374.414
$ws.Cells.Item(62, 10).Value = 22297.166  # J62: 14504.7 -> 22297.166
$ws.Cells.Item(62, 12).Value = 22297.166  # L62: 14504.7 -> 22297.166
$ws.Cells.Item(62, 14).Value = -23545.166  # N62: -15752.7 -> -23545.166
$ws.Cells.Item(65, 8).Value = 69374.414  # H65: 56696.617 -> 69374.414
$ws.Cells.Item(65, 10).Value = 22297.166  # J65: 14504.7 -> 22297.166
$ws.Cells.Item(65, 12).Value = 111485.83  # L65: 72523.5 -> 111485.83
$ws.Cells.Item(65, 14).Value = -117725.83  # N65: -78763.5 -> -117725.83
